# 18 Has Registered User Valid Membership on Date.docx
# "Remove last comments, seperate optional parts"
#
# 1) Drop the empty <w:sdtEndPr/> from the title content control.
# 2) Remove the "[Optional] "date": ..." comment paragraph from the
#    Acceptance Criteria / Response table cell.
# 3) Append a new "Extra - Optional" section at the end of the document
#    describing the optional date-validation requirement and its error
#    response, using a new bulleted list definition (numId 7).

$d = $word.ActiveDocument

# --- 1. Remove the stray <w:sdtEndPr/> on the title content control ---------
$titlePara = $d.Paragraphs.Item(1)
$titleXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="2172BD85" w14:textId="3BF815E0" w:rsidR="00C36337" w:rsidRPr="00614181" w:rsidRDefault="00F66740" w:rsidP="00C36337"><w:pPr><w:pStyle w:val="Title"/><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:sdt><w:sdtPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:lang w:val="en-US"/></w:rPr><w:alias w:val="Title"/><w:tag w:val=""/><w:id w:val="1552501305"/><w:placeholder><w:docPart w:val="21F798DD9A7576419EDE9B53925C317C"/></w:placeholder><w:dataBinding w:prefixMappings="xmlns:ns0='http://purl.org/dc/elements/1.1/' xmlns:ns1='http://schemas.openxmlformats.org/package/2006/metadata/core-properties' " w:xpath="/ns1:coreProperties[1]/ns0:title[1]" w:storeItemID="{6C3C8BC8-F283-45AE-878A-BAB7291924A1}"/><w:text/></w:sdtPr><w:sdtContent><w:r w:rsidR="00E51824"><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:lang w:val="en-US"/></w:rPr><w:t>18 Has Registered User Valid Membership on Date</w:t></w:r></w:sdtContent></w:sdt></w:p>
'@
$titlePara.Range.InsertXML($titleXml)

# --- 2. Remove the "[Optional] date pattern" comment paragraph in the table -
$count = $d.Paragraphs.Count
for ($i = $count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*[Optional]*Date has invalid pattern*") {
        $p.Range.Delete()
        break
    }
}

# --- 3. Append the new "Extra - Optional" section at the end ----------------
$endRng = $d.Content
$endRng.Collapse(0)
$extraXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading2"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="0"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Extra - </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Optional</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr><w:spacing w:line="256" w:lineRule="auto"/><w:rPr><w:lang w:val="en-US" w:eastAsia="nl-NL"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US" w:eastAsia="nl-NL"/></w:rPr><w:t xml:space="preserve">As an optional requirement, you can try to validate the </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US" w:eastAsia="nl-NL"/></w:rPr><w:t>date</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US" w:eastAsia="nl-NL"/></w:rPr><w:t xml:space="preserve"> (search for yourself on how you can validate a date)</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US" w:eastAsia="nl-NL"/></w:rPr><w:t>. The following error</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US" w:eastAsia="nl-NL"/></w:rPr><w:t xml:space="preserve"> is</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US" w:eastAsia="nl-NL"/></w:rPr><w:t xml:space="preserve"> returned when the input is not valid:</w:t></w:r></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="7"/></w:numPr><w:textAlignment w:val="baseline"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:lang w:val="en-US"/></w:rPr><w:t>&#8220;date&#8221;: &#8220;Date has invalid pattern: YYYY-MM-DD</w:t></w:r></w:p>
'@
$endRng.InsertXML($extraXml)

# --- 4. Mark the Heading4 style as a Quick Style (adds <w:qFormat/>) --------
$d.Styles("Heading 4").QuickStyle = $true

Write-Output "done"
